$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New test-case rows appended below the existing "TestCase_B12" row (row 13).
$data = @(
  @("TestCase_B15", "To verify that 10 MORE button is working correctly", "Y", "FAIL"),
  @("TestCase_B16", "To verify that 10 MORE button is not present in search results page if the total search results is less than or equal to 10", "Y", "PASS"),
  @("TestCase_B17", "To verify that 10 MORE button is present in search results page if total search results is more than 10", "Y", "PASS"),
  @("TestCase_B18", "To verify that Times cited,Views,Comments fields are getting displayed for each document in search results page", "Y", "PASS"),
  @("TestCase_B19", "To verify that user is able to expand and collapse SORT BY drop down", "Y", "PASS"),
  @("TestCase_B20", "To verify that user is able to sort the documents by TIMES CITED field", "Y", "PASS"),
  @("TestCase_B21", "To verify that RESET button in the left navigation pane in search results page is working correctly", "Y", "PASS"),
  @("TestCase_B22", "To verify that MORE and LESS links in the left navigation pane are working correctly", "Y", "PASS"),
  @("TestCase_B23", "To verify that user is able to collapse and expand the filters in the left navigation pane in search results page", "Y", "FAIL")
)

$startRow = 14
$r = $startRow
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $r++
}
$endRow = $r - 1

# Match the bordered (no fill) formatting already used on the existing data
# rows (e.g. D13) by copying its format onto the newly added rows.
$src = $ws.Range("D13")
$dst = $ws.Range("A" + $startRow + ":D" + $endRow)
$src.Copy()
$dst.PasteSpecial(-4122)

# Restore the active selection to the new last-used cell.
$ws.Range("C18").Select() | Out-Null
